$d = $word.ActiveDocument

# Replace the full text of a single paragraph (by 1-based paragraph
# index) with new text, while leaving paragraph/run formatting intact.
# Using Range.Text (rather than Find & Replace) also sidesteps the
# AutoCorrect "smart quotes" substitution that would otherwise turn a
# straight apostrophe into a curly one on replacement.
function Set-ParagraphText($index, $oldText, $newText) {
    $para = $d.Paragraphs.Item($index)
    # Paragraph.Range.Text includes the trailing paragraph-mark ("\r"),
    # so compare against the text with that mark appended.
    $current = $para.Range.Text
    if ($current -ne ($oldText + "`r")) {
        throw "Paragraph $index text mismatch: expected '$oldText', found '$current'"
    }
    $para.Range.Text = $newText
}

# 1) "English" -> "Inglese"
#    "English" appears twice in the document: once inside a hyperlink
#    (paragraph 1, which must stay untouched) and once as the
#    standalone section heading (paragraph 3). Only the heading changes.
Set-ParagraphText 3 "English" "Inglese"

# 2) "Don't delay! Book your spot today!" -> Italian
Set-ParagraphText 15 "Don$([char]8217)t delay! Book your spot today!" `
    "Non aspettare! Prenota il tuo posto oggi stesso!"

# 3) "We look forward to seeing you at [EVENT NAME]! " -> Italian
#    This literal sentence (with the placeholder written out as plain
#    text in a single run) is paragraph 38. An earlier, visually
#    identical sentence (paragraph 27) has "[EVENT NAME]" split into
#    its own highlighted run and must remain in English, so it is left
#    alone.
Set-ParagraphText 38 "We look forward to seeing you at [EVENT NAME]! " `
    "Non vediamo l'ora di incontrarti all'evento [EVENT NAME]! "

# 4) "If you have any questions, please contact your country manager:" -> Italian
Set-ParagraphText 39 "If you have any questions, please contact your country manager:" `
    "Per qualsiasi domanda, contatta il tuo country manager:"
